$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 18
$ws1.Range("F5").Value = 15765
$ws1.Range("F9").Value = 15461
$ws1.Range("F11").Value = 9066
$ws1.Range("F16").Value = 203
$ws1.Range("F18").Value = 205
$ws1.Range("F27").Value = 16
$ws1.Range("F36").Value = 327
$ws1.Range("F37").Value = 459

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 71

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 18
$ws4.Range("F5").Value = 15765
$ws4.Range("F9").Value = 15461
$ws4.Range("F11").Value = 9066
$ws4.Range("F16").Value = 203
$ws4.Range("F18").Value = 205
$ws4.Range("F27").Value = 16
$ws4.Range("F32").Value = 71
$ws4.Range("F38").Value = 327
$ws4.Range("F39").Value = 459
